# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handback DateTime"
# timestamp values recorded during handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file
$wsOverview.Range("G2").Value = "2016-10-24 09:49:13"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsZhCn.Range("H2").Value = "2016-10-24 09:49:01"
$wsZhCn.Range("K2").Value = "2016-10-24 09:49:43"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsDeDe.Range("H2").Value = "2016-10-24 09:49:13"
$wsDeDe.Range("K2").Value = "2016-10-24 09:50:00"
